# "discussed with developers about test cases, and executed them"
#
# The "Spint(43) - Day 10" block (rows 57-60 on the "Test Summary" sheet)
# had its Total/Execution/Review counts still blank. After syncing with
# the dev team the numbers were filled in like every earlier day's block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# Day 10 totals that were missing.
$ws.Range("C58").Value = 7070
$ws.Range("C59").Value = 2550
$ws.Range("C60").Value = 2550

# Leave the view scrolled further down / selection where the author
# ended up after entering the numbers.
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M46").Select()
